$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add a comment anchored on "4005823 " (the first occurrence, in the
#    TOBACCO concept-id table cell) saying:
#    "Would concept 4041306 (Tobacco use and exposure) be better?"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("4005823", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)
$d.Comments.Add($rng, "Would concept 4041306 (Tobacco use and exposure) be better?") | Out-Null
$cmt = $d.Comments.Item(1)
$cmt.Author = "Don Torok"
$cmt.Initial = "DT"

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete its old location (end of the
#    SMOKING heading paragraph) and re-add it around "4219336 " in the
#    TOBACCO TYPE concept-id table cell.
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$rng2 = $d.Content
$rng2.Find.Execute("4219336", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.MoveEnd(1, 1)
$d.Bookmarks.Add("_GoBack", $rng2) | Out-Null

Write-Host "done"
